$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern for cells whose new value looks numeric but must stay text:
# format as Text first, assign, then drop back to the Normal style so no
# stray number-format / quote-prefix styling is left behind on the cell.

# Row 4
$ws.Range("D4").Value = "Cyanide Cabinet"
$ws.Range("E4").Value = "E1-3399"
$ws.Range("F4").Value = "2-(DIMETHYLAMINO)PYRIDINE"
$ws.Range("G4").Value = 22

# Row 5
$ws.Range("C5").Value = "13C"
$ws.Range("D5").Value = "13C"
$ws.Range("E5").Value = "E1-3540"
$ws.Range("F5").Value = "2-BROMOBUTANE"

# Row 6
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "2"
$ws.Range("B6").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "203"
$ws.Range("C6").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "203"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "D1-8364"
$ws.Range("F6").Value = "2,6-DICHLOROINDOPHENOL SODIUM DERIVATIVE"

# Row 7
$ws.Range("C7").Value = "13A"
$ws.Range("D7").Value = "13A"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "17040000"
$ws.Range("E7").Style = "Normal"

$ws.Range("F7").Value = "2,6-DICHLOROINDOPHENOL SODIUM DERIVATIVE"

# Row 8
$ws.Range("B8").Value = "Basement"
$ws.Range("C8").Value = "13C"
$ws.Range("D8").Value = "13C"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "16020008"
$ws.Range("E8").Style = "Normal"

$ws.Range("F8").Value = "ACETIC ACID 99+%"
$ws.Range("G8").Value = 33
